# Baza.xlsx — "Create ClassDiagram and databaseDiagram"
#
# The "Instrument_Reservation" legend block (a small 2-column mapping of
# field -> key-type, living in C8:D10) is relocated so it sits alongside
# the other legend block that already occupies column A on rows 12-14
# (PK - primary key / Fk - foreign key / N - nullable). After the move,
# both legends live together on rows 12-14, and rows 8-10 are empty again.
#
# The active worksheet selection also moves from C14 to G16:G17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relocate the C8:D10 legend block down to C12:D14 -------------------
# Use Range.Copy(Destination) (not clipboard Copy/Paste) so the destination
# cell picks up the exact same style index as the source cell instead of
# Excel fabricating a new one.
$ws.Range("C8").Copy($ws.Range("C12"))      # "Instrument_Reservation" header
$ws.Range("C9").Copy($ws.Range("C13"))      # "Instrument.id - FK"
$ws.Range("C10").Copy($ws.Range("C14"))     # "Reservation.id - FK"
$ws.Range("D9:D10").Copy($ws.Range("D13:D14"))  # merged "PK" cell, carries the merge along

# --- 2. Clear the now-vacated source rows (8, 9, 10) ------------------------
$ws.Range("D9:D10").UnMerge()
$ws.Range("C8:D10").Clear()

# Fully reset row-level formatting (e.g. row 8's extra header height) by
# deleting and immediately re-inserting each row, which leaves row numbering
# and every other row untouched but drops the stale row metadata.
$ws.Rows(8).Delete()
$ws.Rows(8).Insert()
$ws.Rows(9).Delete()
$ws.Rows(9).Insert()
$ws.Rows(10).Delete()
$ws.Rows(10).Insert()

# --- 3. Row 12 becomes the new header row for the merged legend -------------
$ws.Rows(12).RowHeight = 15.75

# --- 4. Update the saved selection ------------------------------------------
$ws.Range("G16:G17").Select()
